$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 56.5
$ws.Range("C4").Value = "Madrid"
$ws.Range("D4").Value = 918.97
$ws.Range("E4").Value = "checkin"
$ws.Range("F4").Value = 96.22
